# Update the "Förändrad" (Changed) date column C for rows 2-5
# from serial date 45233 (2023-11-03) to 45243 (2023-11-13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 5; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    $cell.Value = 45243
}
